$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-6 from 2023-11-03 (45233)
# to 2023-11-13 (45243), keeping existing formatting/style.
$newDateSerial = 45243

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDateSerial
}
